$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: headers (Spanish, capitalized, with accents)
$ws.Range("A1").Value = "Comarca nombre"
$ws.Range("B1").Value = "Comarca código"
$ws.Range("C1").Value = "Número habitaciones"
$ws.Range("D1").Value = "Provincia código"
$ws.Range("E1").Value = "Aragón"
$ws.Range("F1").Value = "Municipio código"
$ws.Range("G1").Value = "Provincia nombre"
$ws.Range("H1").Value = "Municipio nombre"

# Row 2: sdmx-dimension:refArea / iaest-dimension:numero-habitaciones / null
$ws.Range("A2").Value = "sdmx-dimension:refArea"
$ws.Range("B2").Value = "null"
$ws.Range("C2").Value = "iaest-dimension:numero-habitaciones"
$ws.Range("D2").Value = "null"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "null"
$ws.Range("G2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "sdmx-dimension:refArea"

# Row 3: dim / null
$ws.Range("A3").Value = "dim"
$ws.Range("B3").Value = "null"
$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "dim"
$ws.Range("F3").Value = "null"
$ws.Range("G3").Value = "dim"
$ws.Range("H3").Value = "dim"

# Row 4: URI-* / skos:Concept / null
$ws.Range("A4").Value = "URI-comarca"
$ws.Range("B4").Value = "null"
$ws.Range("C4").Value = "skos:Concept"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("F4").Value = "null"
$ws.Range("G4").Value = "URI-Provincia"
$ws.Range("H4").Value = "URI-Municipio"

# Row 5: mapping file reference moved from F5 to C5
# Copy F5's formatting (style s="1") onto C5 before clearing F5 and writing the value
$ws.Range("F5").Copy()
$ws.Range("C5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C5").Value = "mapping-numero-habitaciones.xlsx"
$ws.Range("F5").ClearContents()
